$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2026-01-05 Monday" "2026-01-06 Tuesday"

Replace-Text "290×2=580" "803×9=7227"
Replace-Text "111×6=666" "953×7=6671"
Replace-Text "839×7=5873" "270×4=1080"
Replace-Text "722×9=6498" "823×3=2469"
Replace-Text "465×8=3720" "349×8=2792"

Replace-Text "743×6=4458" "486×6=2916"
Replace-Text "968×3=2904" "932×5=4660"
Replace-Text "794×9=7146" "104×7=728"
Replace-Text "628×3=1884" "584×7=4088"
Replace-Text "146×6=876" "453×5=2265"

Replace-Text "481×4=1924" "270×8=2160"
Replace-Text "211×7=1477" "204×8=1632"
Replace-Text "414×6=2484" "747×5=3735"
Replace-Text "278×9=2502" "166×6=996"
Replace-Text "785×3=2355" "293×8=2344"

Replace-Text "766×9=6894" "974×5=4870"
Replace-Text "459×8=3672" "734×7=5138"
Replace-Text "120×5=600" "926×9=8334"
Replace-Text "182×3=546" "526×2=1052"
Replace-Text "118×2=236" "493×9=4437"

Replace-Text "309×4=1236" "531×7=3717"
Replace-Text "228×5=1140" "951×7=6657"
Replace-Text "376×6=2256" "568×9=5112"
Replace-Text "281×7=1967" "934×5=4670"
Replace-Text "329×7=2303" "131×5=655"
